# The commit swaps the data of row 2 and row 3 (the two records change
# places in the sheet), so every column that differs between the two
# records needs its row-2 and row-3 values exchanged.
#
# Columns C, D, T, U, V, W, Z, AB, AD, AE, AG, I, AT, AY are identical
# between the two rows (or blank in both), so they do not need to be
# touched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 gets what used to be row 3's data ---
$ws.Range("A2").Value = 105312389
$ws.Range("B2").Value = 78503
$ws.Range("E2").Value = 6456
$ws.Range("F2").Value = "Skinnlav"
$ws.Range("G2").Value = "Leptogium saturninum"
$ws.Range("H2").Value = "(Dicks.) Nyl."
$ws.Range("P2").Value = "Storrönningen, Hls"
$ws.Range("Q2").Value = 616060.6447056353
$ws.Range("R2").Value = 6863194.20718522
$ws.Range("S2").Value = 5

$ws.Range("Y2").NumberFormat = "@"
$ws.Range("Y2").Value = "2022-06-03"
$ws.Range("AA2").NumberFormat = "@"
$ws.Range("AA2").Value = "2022-06-03"

$ws.Range("AS2").Value = "Henrik Tykosson"
$ws.Range("AW2").Value = "Helene Andersson"
$ws.Range("AX2").Value = "Henrik Tykosson"

# --- Row 3 gets what used to be row 2's data ---
$ws.Range("A3").Value = 86851042
$ws.Range("B3").Value = 96251
$ws.Range("E3").Value = 219790
$ws.Range("F3").Value = "Fläcknycklar"
$ws.Range("G3").Value = "Dactylorhiza maculata"
$ws.Range("H3").Value = "(L.) Soó"
$ws.Range("P3").Value = "Längs E4 mellan Kongberget och Gnarp, Hls"
$ws.Range("Q3").Value = 615689.084506036
$ws.Range("R3").Value = 6862637.86594828
$ws.Range("S3").Value = 10

$ws.Range("Y3").NumberFormat = "@"
$ws.Range("Y3").Value = "2018-06-27"
$ws.Range("AA3").NumberFormat = "@"
$ws.Range("AA3").Value = "2018-06-27"

$ws.Range("AS3").Value = ""
$ws.Range("AW3").Value = "Sofia Lundman"
$ws.Range("AX3").Value = "Sofia Lundman, Oskar Wallströmer"
